# Auto-generated edit script: updates odds values per the commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 2
$ws.Range("G2").Value = 1.83
$ws.Range("H2").Value = 3.2
$ws.Range("I2").Value = 4.75
$ws.Range("J2").Value = 2.6
$ws.Range("L2").Value = 5.5
$ws.Range("M2").Value = 1.1
$ws.Range("N2").Value = 7
$ws.Range("Q2").Value = 2.6
$ws.Range("R2").Value = 1.48
$ws.Range("U2").Value = 2.25
$ws.Range("V2").Value = 1.57
$ws.Range("X2").Value = 7.5
$ws.Range("Z2").Value = 15
$ws.Range("AE2").Value = 21
$ws.Range("AG2").Value = 9.5
$ws.Range("AH2").Value = 21
$ws.Range("AI2").Value = 17
$ws.Range("AK2").Value = 41
$ws.Range("AN2").Value = 3.6
$ws.Range("AO2").Value = 11
$ws.Range("AQ2").Value = 41
$ws.Range("AU2").Value = 9.5
$ws.Range("AX2").Value = 29
$ws.Range("AZ2").Value = 101
$ws.Range("BA2").Value = 151

# Row 3
$ws.Range("G3").Value = 1.95
$ws.Range("H3").Value = 3.2
$ws.Range("I3").Value = 4.2
$ws.Range("J3").Value = 2.75
$ws.Range("L3").Value = 4.75
$ws.Range("U3").Value = 2.1
$ws.Range("V3").Value = 1.67
$ws.Range("X3").Value = 8
$ws.Range("Y3").Value = 9.5
$ws.Range("Z3").Value = 17
$ws.Range("AA3").Value = 19
$ws.Range("AB3").Value = 34
$ws.Range("AE3").Value = 19
$ws.Range("AF3").Value = 67
$ws.Range("AG3").Value = 9.5
$ws.Range("AH3").Value = 19
$ws.Range("AI3").Value = 15
$ws.Range("AJ3").Value = 41
$ws.Range("AN3").Value = 3.75
$ws.Range("AO3").Value = 11
$ws.Range("AQ3").Value = 41
$ws.Range("AU3").Value = 9
$ws.Range("AV3").Value = 67
$ws.Range("AW3").Value = 6
$ws.Range("AX3").Value = 26
$ws.Range("AZ3").Value = 81
$ws.Range("BA3").Value = 126

# Row 5
$ws.Range("Q5").Value = 2.25
$ws.Range("R5").Value = 1.62

# Row 6
$ws.Range("M6").Value = 1.06
$ws.Range("N6").Value = 8.380000000000001
$ws.Range("O6").Value = 1.32
$ws.Range("P6").Value = 2.85
$ws.Range("Q6").Value = 1.93
$ws.Range("R6").Value = 1.7
$ws.Range("U6").Value = 1.72
$ws.Range("V6").Value = 1.88
$ws.Range("W6").Value = 7.3
$ws.Range("X6").Value = 10.25
$ws.Range("AC6").Value = 8.75
$ws.Range("AE6").Value = 14
$ws.Range("AG6").Value = 9.25
$ws.Range("AH6").Value = 17
$ws.Range("AK6").Value = 30
$ws.Range("AL6").Value = 37
$ws.Range("AN6").Value = 4.05
$ws.Range("AP6").Value = 18
$ws.Range("AR6").Value = 70
$ws.Range("AT6").Value = 2.57
$ws.Range("AU6").Value = 6.7
$ws.Range("AY6").Value = 23
$ws.Range("AZ6").Value = 90
$ws.Range("BB6").Value = 300

# Row 7
$ws.Range("P7").Value = 4.25

# Row 10
$ws.Range("Q10").Value = 2.15
$ws.Range("R10").Value = 1.67

# Row 11
$ws.Range("G11").Value = 1.36
$ws.Range("H11").Value = 4.5
$ws.Range("I11").Value = 7.5
$ws.Range("J11").Value = 1.91
$ws.Range("L11").Value = 8
$ws.Range("U11").Value = 2.38
$ws.Range("V11").Value = 1.53
$ws.Range("W11").Value = 5.5
$ws.Range("X11").Value = 5.5
$ws.Range("Z11").Value = 8.5
$ws.Range("AD11").Value = 9.5
$ws.Range("AE11").Value = 26
$ws.Range("AF11").Value = 101
$ws.Range("AK11").Value = 67
$ws.Range("AO11").Value = 6.5
$ws.Range("AS11").Value = 201
$ws.Range("AV11").Value = 81
$ws.Range("AW11").Value = 9
$ws.Range("BA11").Value = 251

# Row 12
$ws.Range("M12").Value = 1.06
$ws.Range("N12").Value = 10
$ws.Range("Q12").Value = 2.08
$ws.Range("R12").Value = 1.73

# Row 16
$ws.Range("O16").Value = 1.33
$ws.Range("P16").Value = 3.25

# Row 17
$ws.Range("M17").Value = 1.03
$ws.Range("N17").Value = 15

# Row 19
$ws.Range("Q19").Value = 2
$ws.Range("R19").Value = 1.85

# Row 21
$ws.Range("Q21").Value = 2.5
$ws.Range("R21").Value = 1.5

# Row 27
$ws.Range("G27").Value = 2.75
$ws.Range("J27").Value = 3.55
$ws.Range("K27").Value = 1.83
$ws.Range("L27").Value = 3.4
$ws.Range("O27").Value = 1.5
$ws.Range("P27").Value = 2.25
$ws.Range("Q27").Value = 2.45
$ws.Range("S27").Value = 1.55
$ws.Range("T27").Value = 2.15
$ws.Range("W27").Value = 6.3
$ws.Range("X27").Value = 12
$ws.Range("Y27").Value = 11
$ws.Range("AA27").Value = 32
$ws.Range("AB27").Value = 50
$ws.Range("AC27").Value = 6
$ws.Range("AG27").Value = 6.8
$ws.Range("AH27").Value = 13
$ws.Range("AI27").Value = 10.5
$ws.Range("AK27").Value = 28
$ws.Range("AL27").Value = 45
$ws.Range("AO27").Value = 16.5
$ws.Range("AP27").Value = 29
$ws.Range("AQ27").Value = 90
$ws.Range("AS27").Value = 500
$ws.Range("AU27").Value = 7.5
$ws.Range("AW27").Value = 4.4
$ws.Range("AX27").Value = 15.5
$ws.Range("AY27").Value = 26
$ws.Range("AZ27").Value = 75
$ws.Range("BA27").Value = 120
$ws.Range("BB27").Value = 400
